# Applies the "Automatic update of files" change:
#  - Column C ("Förändrad") is bumped from 46072 to 46073 for every data row (2-14).
#  - Rows 6-12 and 14 are re-shuffled: the full record (Beteckning/A, Datum/B, Area/G)
#    that used to live in one row now lives in a different row. Row 13 is unchanged.
#
# Rather than trying to move ranges around, we just write the final values for the
# cells that differ, directly, row by row. (Note: this runtime's PowerShell does not
# support named parameters, so Set-Row uses positional parameters; pass $null to
# skip a column.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        $A,
        $B,
        $C,
        $G
    )
    if ($null -ne $A) { $ws.Cells.Item($Row, 1).Value = $A }
    if ($null -ne $B) { $ws.Cells.Item($Row, 2).Value = $B }
    if ($null -ne $C) { $ws.Cells.Item($Row, 3).Value = $C }
    if ($null -ne $G) { $ws.Cells.Item($Row, 7).Value = $G }
}

# Rows 2-5: only column C changes (46072 -> 46073).
Set-Row 2 $null $null 46073 $null
Set-Row 3 $null $null 46073 $null
Set-Row 4 $null $null 46073 $null
Set-Row 5 $null $null 46073 $null

# Rows 6-12 and 14 are reshuffled (A, B, G move together), plus C bumps everywhere.
Set-Row 6  "A 27724-2022" 44743.48386574074 46073 1.3
Set-Row 7  "A 46779-2025" 45926              46073 1.5
Set-Row 8  "A 56917-2025" 45978.58453703704 46073 0.7
Set-Row 9  "A 56948-2025" 45978.64356481482 46073 4.7
Set-Row 10 "A 64445-2023" 45280              46073 3.7
Set-Row 11 "A 31120-2023" 45113              46073 0.2
Set-Row 12 "A 64431-2023" 45280              46073 0.5

# Row 13: only column C changes.
Set-Row 13 $null $null 46073 $null

# Row 14: gets the record that used to be on row 9.
Set-Row 14 "A 50934-2024" 45602 46073 0.6
